# Rename the worksheet: "data" -> "data_example"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "data_example"

# The existing header cells (B1:D1) carried a custom bold/white-on-blue
# style; the new version of the sheet drops that formatting entirely so
# the headers render with the workbook's default style.
$ws.Range("B1:D1").ClearFormats()

# Two brand-new header columns are appended after Pressure.
$ws.Range("E1").Value = "mitame"
$ws.Range("F1").Value = "sigma"

# Replace/extend the single sample data row with the full 8-row table
# (columns A..F): label, Temperature, Time, Pressure, mitame, sigma.
$rows = @(
    @("d_1",  270, 1,  2, 4, 1),
    @("c_02", 270, 1,  4, 3, 0.5),
    @("c_11", 270, 13, 4, 3, 0.2),
    @("c_21", 280, 13, 2, 5, 0.4),
    @("c_23", 280, 13, 4, 2, 0.3),
    @("c_33", 290, 13, 2, 5, 1.1),
    @("c_75", 330, 5,  2, 4, 0.2),
    @("c_86", 340, 1,  4, 1, 0.1)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
